$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-08-18 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-19 Saturday", 2) | Out-Null

# Update each table cell value (5 columns x 20 rows)
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "49-11=38"
$t.Cell(1,2).Range.Text = "21-1=20"
$t.Cell(1,3).Range.Text = "92-7=85"
$t.Cell(1,4).Range.Text = "51+4=55"
$t.Cell(1,5).Range.Text = "78-27=51"
$t.Cell(2,1).Range.Text = "75+13=88"
$t.Cell(2,2).Range.Text = "0+44=44"
$t.Cell(2,3).Range.Text = "26+13=39"
$t.Cell(2,4).Range.Text = "46+3=49"
$t.Cell(2,5).Range.Text = "93-55=38"
$t.Cell(3,1).Range.Text = "83-27=56"
$t.Cell(3,2).Range.Text = "68-48=20"
$t.Cell(3,3).Range.Text = "85-17=68"
$t.Cell(3,4).Range.Text = "34+50=84"
$t.Cell(3,5).Range.Text = "88-80=8"
$t.Cell(4,1).Range.Text = "21+9=30"
$t.Cell(4,2).Range.Text = "83-8=75"
$t.Cell(4,3).Range.Text = "46-35=11"
$t.Cell(4,4).Range.Text = "0+52=52"
$t.Cell(4,5).Range.Text = "46+24=70"
$t.Cell(5,1).Range.Text = "66-64=2"
$t.Cell(5,2).Range.Text = "64-7=57"
$t.Cell(5,3).Range.Text = "3+32=35"
$t.Cell(5,4).Range.Text = "89-63=26"
$t.Cell(5,5).Range.Text = "62+12=74"
$t.Cell(6,1).Range.Text = "4+59=63"
$t.Cell(6,2).Range.Text = "62+8=70"
$t.Cell(6,3).Range.Text = "78+13=91"
$t.Cell(6,4).Range.Text = "48-17=31"
$t.Cell(6,5).Range.Text = "36-19=17"
$t.Cell(7,1).Range.Text = "86-46=40"
$t.Cell(7,2).Range.Text = "0+57=57"
$t.Cell(7,3).Range.Text = "14+14=28"
$t.Cell(7,4).Range.Text = "14-7=7"
$t.Cell(7,5).Range.Text = "8+45=53"
$t.Cell(8,1).Range.Text = "67-23=44"
$t.Cell(8,2).Range.Text = "77-73=4"
$t.Cell(8,3).Range.Text = "77-24=53"
$t.Cell(8,4).Range.Text = "72-65=7"
$t.Cell(8,5).Range.Text = "18+14=32"
$t.Cell(9,1).Range.Text = "16+59=75"
$t.Cell(9,2).Range.Text = "49-42=7"
$t.Cell(9,3).Range.Text = "81-72=9"
$t.Cell(9,4).Range.Text = "72+23=95"
$t.Cell(9,5).Range.Text = "37+23=60"
$t.Cell(10,1).Range.Text = "78-70=8"
$t.Cell(10,2).Range.Text = "16+41=57"
$t.Cell(10,3).Range.Text = "78-52=26"
$t.Cell(10,4).Range.Text = "36-9=27"
$t.Cell(10,5).Range.Text = "53+13=66"
$t.Cell(11,1).Range.Text = "35+56=91"
$t.Cell(11,2).Range.Text = "8+67=75"
$t.Cell(11,3).Range.Text = "52-28=24"
$t.Cell(11,4).Range.Text = "18+22=40"
$t.Cell(11,5).Range.Text = "40+51=91"
$t.Cell(12,1).Range.Text = "5+49=54"
$t.Cell(12,2).Range.Text = "24+72=96"
$t.Cell(12,3).Range.Text = "6+84=90"
$t.Cell(12,4).Range.Text = "84-25=59"
$t.Cell(12,5).Range.Text = "47-39=8"
$t.Cell(13,1).Range.Text = "55-54=1"
$t.Cell(13,2).Range.Text = "9+78=87"
$t.Cell(13,3).Range.Text = "5+91=96"
$t.Cell(13,4).Range.Text = "87-21=66"
$t.Cell(13,5).Range.Text = "87+2=89"
$t.Cell(14,1).Range.Text = "94-67=27"
$t.Cell(14,2).Range.Text = "14+6=20"
$t.Cell(14,3).Range.Text = "15+39=54"
$t.Cell(14,4).Range.Text = "2+49=51"
$t.Cell(14,5).Range.Text = "16+60=76"
$t.Cell(15,1).Range.Text = "75-31=44"
$t.Cell(15,2).Range.Text = "31+67=98"
$t.Cell(15,3).Range.Text = "12+70=82"
$t.Cell(15,4).Range.Text = "73-39=34"
$t.Cell(15,5).Range.Text = "16-9=7"
$t.Cell(16,1).Range.Text = "37+15=52"
$t.Cell(16,2).Range.Text = "9+43=52"
$t.Cell(16,3).Range.Text = "15-5=10"
$t.Cell(16,4).Range.Text = "45+15=60"
$t.Cell(16,5).Range.Text = "0+16=16"
$t.Cell(17,1).Range.Text = "46-37=9"
$t.Cell(17,2).Range.Text = "21+55=76"
$t.Cell(17,3).Range.Text = "6+25=31"
$t.Cell(17,4).Range.Text = "48-2=46"
$t.Cell(17,5).Range.Text = "87-32=55"
$t.Cell(18,1).Range.Text = "75-21=54"
$t.Cell(18,2).Range.Text = "97-49=48"
$t.Cell(18,3).Range.Text = "70+13=83"
$t.Cell(18,4).Range.Text = "84-49=35"
$t.Cell(18,5).Range.Text = "33+48=81"
$t.Cell(19,1).Range.Text = "65-38=27"
$t.Cell(19,2).Range.Text = "13+19=32"
$t.Cell(19,3).Range.Text = "9+22=31"
$t.Cell(19,4).Range.Text = "76-0=76"
$t.Cell(19,5).Range.Text = "24+21=45"
$t.Cell(20,1).Range.Text = "66-18=48"
$t.Cell(20,2).Range.Text = "80-16=64"
$t.Cell(20,3).Range.Text = "20+30=50"
$t.Cell(20,4).Range.Text = "95+4=99"
$t.Cell(20,5).Range.Text = "25+1=26"
